$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98
$prev = $row - 1

# Copy the formatting (styles) of the row above into the new row so
# that the bordered/bold "Indice" column and the date-formatted
# "data_partida" column keep their existing style indices.
$ws.Range("A$prev`:V$prev").Copy()
$ws.Range("A$row`:V$row").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = 97
$ws.Cells.Item($row, 2).Value = "belgium"
$ws.Cells.Item($row, 3).Value = "jupiler-pro-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"

$ws.Cells.Item($row, 5).Value = 45233.86458333334

$ws.Cells.Item($row, 6).Value = "St. Truiden"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Eupen"
$ws.Cells.Item($row, 9).Value = 1

$ws.Cells.Item($row, 10).Value = 1.66
$ws.Cells.Item($row, 11).Value = "28/10/2023 18:42"
$ws.Cells.Item($row, 12).Value = 1.67
$ws.Cells.Item($row, 13).Value = "03/11/2023 20:39"

$ws.Cells.Item($row, 14).Value = 4.11
$ws.Cells.Item($row, 15).Value = "28/10/2023 18:42"
$ws.Cells.Item($row, 16).Value = 4.07
$ws.Cells.Item($row, 17).Value = "03/11/2023 20:39"

$ws.Cells.Item($row, 18).Value = 5
$ws.Cells.Item($row, 19).Value = "28/10/2023 18:42"
$ws.Cells.Item($row, 20).Value = 5.17
$ws.Cells.Item($row, 21).Value = "03/11/2023 20:39"

$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/st-truiden-eupen/UsaIchee/"
